$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(10, "aa", "dssdfd"),
    @(11, "bb", "gfdsgdfs"),
    @(12, "cc", "sfdgf"),
    @(13, "dd", "sdfgfsd"),
    @(14, "ee", "sdfgsfd"),
    @(15, "ff", "sfgsfg"),
    @(16, "gg", "fsdgdfs"),
    @(17, "hh", "sfgfsdfg"),
    @(18, "ii", "sdfg"),
    @(19, "jj", "sfgffd")
)

$row = 11
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
